# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
#
# Rows 285-293 on "Brazil Serie B" have their match-data columns (B:AC)
# permuted between rows while the leading index column (A) stays put.
# Snapshot every source row first (Value2, to avoid the shared-string /
# reflection quirks of the plain .Value getter), then write them all back
# so the cyclic permutation doesn't clobber a row before it's been read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row285 = $ws.Range("B285:AC285").Value2
$row286 = $ws.Range("B286:AC286").Value2
$row287 = $ws.Range("B287:AC287").Value2
$row288 = $ws.Range("B288:AC288").Value2
$row289 = $ws.Range("B289:AC289").Value2
$row290 = $ws.Range("B290:AC290").Value2
$row292 = $ws.Range("B292:AC292").Value2
$row293 = $ws.Range("B293:AC293").Value2

# row 291 is unchanged (maps to itself), so it is left untouched.

$ws.Range("B285:AC285").Value2 = $row286
$ws.Range("B286:AC286").Value2 = $row285
$ws.Range("B287:AC287").Value2 = $row293
$ws.Range("B288:AC288").Value2 = $row287
$ws.Range("B289:AC289").Value2 = $row288
$ws.Range("B290:AC290").Value2 = $row292
$ws.Range("B292:AC292").Value2 = $row290
$ws.Range("B293:AC293").Value2 = $row289
